$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 538, shifting existing rows 538:653 down to 539:654.
$ws.Rows.Item(538).Insert()

# Populate the newly inserted row 538 with the new data record.
$ws.Range("A538").Value = 5
$ws.Range("B538").Value = "Macroferia Regional de Talca"
$ws.Range("C538").Value = "Maule"
$ws.Range("D538").Value = 45173
$ws.Range("E538").Value = 7
$ws.Range("F538").Value = 100114014
$ws.Range("G538").Value = "Betarraga"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "Primera"
$ws.Range("J538").Value = 6000
$ws.Range("K538").Value = 450
$ws.Range("L538").Value = 500
$ws.Range("M538").Value = 475
$ws.Range("N538").Value = "$/paquete 5 unidades"
$ws.Range("O538").Value = "Región del Maule"
$ws.Range("P538").Value = 95
$ws.Range("Q538").Value = 5
$ws.Range("R538").Value = "Hortaliza"
